# Adds QC Nuclear decommission policy comments describing the
# Gentilly-2 shutdown (2012) and the blank 2030+ target years that
# allow for a potential refurbishment of Gentilly-2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("X3").Value = "Gentilly-2 shutdown in 2012"
$ws.Range("X4").Value = "2030 and later blank to allow refurbishment of Gentilly-2"
